$d = $word.ActiveDocument

# --- Locate and split the target paragraph text -----------------------
$old = "Las rotaciones se prueban en ambas direcciones (izquierda y derecha), desde 1 hasta 7 bits. No se consideran rotaciones de 8 bits o superiores, ya que una rotación de 8 devuelve el mismo byte original, y valores mayores equivalen a rotaciones más pequeñas (por ejemplo, rotar 9 bits equivale a rotar 1). "
$part1 = "Las rotaciones se prueban en ambas direcciones (izquierda y derecha), desde 1 hasta 7 bits. No se consideran rotaciones de 8 bits o superiores, ya que una rotación de 8 devuelve el mismo byte original"

$r = $d.Content
$found = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $part1, 2)

# $r now spans exactly the replaced "part1" text (same run/format as the
# original sentence). Collapse to its end so we can append the rest of
# the new material, run by run, exactly as described in the diff.
$r.Collapse(0)

# Run 2: " y en el documento dice que" (same formatting as run 1)
$r.InsertAfter(" y en el documento dice que")
$r.Collapse(0)

# Run 3: ": ”" (same formatting)
$r.InsertAfter(": ”")
$r.Collapse(0)

# Run 4: a lone space, distinct run (no explicit Times New Roman/size
# override in the source diff - only the language mark remains).
$spaceStart = $r.End
$r.InsertAfter(" ")
$r.Collapse(0)
$spaceRange = $d.Range($spaceStart, $r.End)
$spaceRange.LanguageID = 14698

# Run 5: "E"
$r.InsertAfter("E")
$r.Collapse(0)

# Run 6: "l máximo número de bits a rotar o desplazar es de 8."
$r.InsertAfter("l máximo número de bits a rotar o desplazar es de 8.")
$r.Collapse(0)

# Run 7: closing quote "”"
$r.InsertAfter("”")
$r.Collapse(0)
